$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 63.18739966666666
$ws.Range("H2").Value = 189.562199
$ws.Range("I2").Value = 0.09596345243430386
$ws.Range("J2").Value = 0.09988075390087989
$ws.Range("M2").Value = 12.431794
$ws.Range("N2").Value = 37.295382
$ws.Range("O2").Value = 0.6267040910788743
$ws.Range("P2").Value = 0.7075740515758999
$ws.Range("Q2").Value = 785.5327360516687
$ws.Range("R2").Value = 7069.794624465018
$ws.Range("S2").Value = 0.06014068823463119
$ws.Range("T2").Value = 0.07067302971210095

$ws.Range("G3").Value = 63.18739966666666
$ws.Range("H3").Value = 189.562199
$ws.Range("I3").Value = 0.09596345243430386
$ws.Range("J3").Value = 0.09988075390087989
$ws.Range("O3").Value = 0.0264162940991436
$ws.Range("P3").Value = 0.0298250554119953
$ws.Range("Q3").Value = 33.11110311139522
$ws.Range("R3").Value = 297.9999280025569
$ws.Range("S3").Value = 0.002534998782273748
$ws.Range("T3").Value = 0.002978949019685608

$ws.Range("G4").Value = 63.18739966666666
$ws.Range("H4").Value = 189.562199
$ws.Range("I4").Value = 0.09596345243430386
$ws.Range("J4").Value = 0.09988075390087989
$ws.Range("M4").Value = 0.03915333333333333
$ws.Range("N4").Value = 0.11746
$ws.Range("O4").Value = 0.001973774193762771
$ws.Range("P4").Value = 0.002228470219130754
$ws.Range("Q4").Value = 2.473997321615555
$ws.Range("R4").Value = 22.26597589454
$ws.Range("S4").Value = 0.0001894101859592101
$ws.Range("T4").Value = 0.0002225812855324387

$ws.Range("G5").Value = 63.18739966666666
$ws.Range("H5").Value = 189.562199
$ws.Range("I5").Value = 0.09596345243430386
$ws.Range("J5").Value = 0.09988075390087989
$ws.Range("M5").Value = 6.8015495
$ws.Range("N5").Value = 13.603099
$ws.Range("O5").Value = 0.3428756056708687
$ws.Range("P5").Value = 0.2580802061075034
$ws.Range("Q5").Value = 429.7722266091168
$ws.Range("R5").Value = 2578.633359654701
$ws.Range("S5").Value = 0.03290352687567953
$ws.Range("T5").Value = 0.0257772455529119

$ws.Range("G6").Value = 63.18739966666666
$ws.Range("H6").Value = 189.562199
$ws.Range("I6").Value = 0.09596345243430386
$ws.Range("J6").Value = 0.09988075390087989
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04027333333333333
$ws.Range("N6").Value = 0.12082
$ws.Range("O6").Value = 0.002030234957350741
$ws.Range("P6").Value = 0.002292216685470609
$ws.Range("Q6").Value = 2.544767209242222
$ws.Range("R6").Value = 22.90290488318
$ws.Range("S6").Value = 0.0001948283557601887
$ws.Range("T6").Value = 0.0002289483306489805

$ws.Range("I7").Value = 0.3063997713314046
$ws.Range("J7").Value = 0.3189072441572365
$ws.Range("M7").Value = 12.431794
$ws.Range("N7").Value = 37.295382
$ws.Range("O7").Value = 0.6267040910788743
$ws.Range("P7").Value = 0.7075740515758999
$ws.Range("Q7").Value = 2508.111625770625
$ws.Range("R7").Value = 22573.00463193563
$ws.Range("S7").Value = 0.1920219901990228
$ws.Range("T7").Value = 0.2256504908252406

$ws.Range("I8").Value = 0.3063997713314046
$ws.Range("J8").Value = 0.3189072441572365
$ws.Range("O8").Value = 0.0264162940991436
$ws.Range("P8").Value = 0.0298250554119953
$ws.Range("S8").Value = 0.008093946471400732
$ws.Range("T8").Value = 0.009511426228276295

$ws.Range("I9").Value = 0.3063997713314046
$ws.Range("J9").Value = 0.3189072441572365
$ws.Range("M9").Value = 0.03915333333333333
$ws.Range("N9").Value = 0.11746
$ws.Range("O9").Value = 0.001973774193762771
$ws.Range("P9").Value = 0.002228470219130754
$ws.Range("Q9").Value = 7.899176138295555
$ws.Range("R9").Value = 71.09258524466
$ws.Range("S9").Value = 0.0006047639616287404
$ws.Range("T9").Value = 0.0007106752962694618

$ws.Range("I10").Value = 0.3063997713314046
$ws.Range("J10").Value = 0.3189072441572365
$ws.Range("M10").Value = 6.8015495
$ws.Range("N10").Value = 13.603099
$ws.Range("O10").Value = 0.3428756056708687
$ws.Range("P10").Value = 0.2580802061075034
$ws.Range("Q10").Value = 1372.21107220763
$ws.Range("R10").Value = 8233.266433245779
$ws.Range("S10").Value = 0.105057007172671
$ws.Range("T10").Value = 0.08230364730127551

$ws.Range("I11").Value = 0.3063997713314046
$ws.Range("J11").Value = 0.3189072441572365
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04027333333333333
$ws.Range("N11").Value = 0.12082
$ws.Range("O11").Value = 0.002030234957350741
$ws.Range("P11").Value = 0.002292216685470609
$ws.Range("Q11").Value = 8.125135884802223
$ws.Range("R11").Value = 73.12622296322
$ws.Range("S11").Value = 0.000622063526681291
$ws.Range("T11").Value = 0.000731004506174667

$ws.Range("G12").Value = 170.2928416666667
$ws.Range("H12").Value = 510.878525
$ws.Range("I12").Value = 0.2586257560429799
$ws.Range("J12").Value = 0.2691830570543736
$ws.Range("M12").Value = 12.431794
$ws.Range("N12").Value = 37.295382
$ws.Range("O12").Value = 0.6267040910788743
$ws.Range("P12").Value = 0.7075740515758999
$ws.Range("Q12").Value = 2117.045527274617
$ws.Range("R12").Value = 19053.40974547155
$ws.Range("S12").Value = 0.1620818193705024
$ws.Range("T12").Value = 0.1904669462955497

$ws.Range("G13").Value = 170.2928416666667
$ws.Range("H13").Value = 510.878525
$ws.Range("I13").Value = 0.2586257560429799
$ws.Range("J13").Value = 0.2691830570543736
$ws.Range("O13").Value = 0.0264162940991436
$ws.Range("P13").Value = 0.0298250554119953
$ws.Range("Q13").Value = 89.23588989739721
$ws.Range("R13").Value = 803.1230090765749
$ws.Range("S13").Value = 0.006831934033244723
$ws.Range("T13").Value = 0.008028399592616985

$ws.Range("G14").Value = 170.2928416666667
$ws.Range("H14").Value = 510.878525
$ws.Range("I14").Value = 0.2586257560429799
$ws.Range("J14").Value = 0.2691830570543736
$ws.Range("M14").Value = 0.03915333333333333
$ws.Range("N14").Value = 0.11746
$ws.Range("O14").Value = 0.001973774193762771
$ws.Range("P14").Value = 0.002228470219130754
$ws.Range("Q14").Value = 6.667532394055556
$ws.Range("R14").Value = 60.0077915465
$ws.Range("S14").Value = 0.0005104688431200198
$ws.Range("T14").Value = 0.0005998664261402461

$ws.Range("G15").Value = 170.2928416666667
$ws.Range("H15").Value = 510.878525
$ws.Range("I15").Value = 0.2586257560429799
$ws.Range("J15").Value = 0.2691830570543736
$ws.Range("M15").Value = 6.8015495
$ws.Range("N15").Value = 13.603099
$ws.Range("O15").Value = 0.3428756056708687
$ws.Range("P15").Value = 0.2580802061075034
$ws.Range("Q15").Value = 1158.255192091496
$ws.Range("R15").Value = 6949.531152548975
$ws.Range("S15").Value = 0.08867646274532306
$ws.Range("T15").Value = 0.06947081884524058

$ws.Range("G16").Value = 170.2928416666667
$ws.Range("H16").Value = 510.878525
$ws.Range("I16").Value = 0.2586257560429799
$ws.Range("J16").Value = 0.2691830570543736
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.04027333333333333
$ws.Range("N16").Value = 0.12082
$ws.Range("O16").Value = 0.002030234957350741
$ws.Range("P16").Value = 0.002292216685470609
$ws.Range("Q16").Value = 6.858260376722223
$ws.Range("R16").Value = 61.7243433905
$ws.Range("S16").Value = 0.0005250710507897225
$ws.Range("T16").Value = 0.0006170258948260221

$ws.Range("G17").Value = 77.473122
$ws.Range("H17").Value = 154.946244
$ws.Range("I17").Value = 0.1176593481802354
$ws.Range("J17").Value = 0.08164152846121862
$ws.Range("M17").Value = 12.431794
$ws.Range("N17").Value = 37.295382
$ws.Range("O17").Value = 0.6267040910788743
$ws.Range("P17").Value = 0.7075740515758999
$ws.Range("Q17").Value = 963.1298932408682
$ws.Range("R17").Value = 5778.779359445209
$ws.Range("S17").Value = 0.07373759485822719
$ws.Range("T17").Value = 0.05776742707015361

$ws.Range("G18").Value = 77.473122
$ws.Range("H18").Value = 154.946244
$ws.Range("I18").Value = 0.1176593481802354
$ws.Range("J18").Value = 0.08164152846121862
$ws.Range("O18").Value = 0.0264162940991436
$ws.Range("P18").Value = 0.0298250554119953
$ws.Range("Q18").Value = 40.59702637608199
$ws.Range("R18").Value = 243.582158256492
$ws.Range("S18").Value = 0.003108123945042634
$ws.Range("T18").Value = 0.002434963110275837

$ws.Range("G19").Value = 77.473122
$ws.Range("H19").Value = 154.946244
$ws.Range("I19").Value = 0.1176593481802354
$ws.Range("J19").Value = 0.08164152846121862
$ws.Range("M19").Value = 0.03915333333333333
$ws.Range("N19").Value = 0.11746
$ws.Range("O19").Value = 0.001973774193762771
$ws.Range("P19").Value = 0.002228470219130754
$ws.Range("Q19").Value = 3.03333097004
$ws.Range("R19").Value = 18.19998582024
$ws.Range("S19").Value = 0.0002322329850930972
$ws.Range("T19").Value = 0.0001819357148201416

$ws.Range("G20").Value = 77.473122
$ws.Range("H20").Value = 154.946244
$ws.Range("I20").Value = 0.1176593481802354
$ws.Range("J20").Value = 0.08164152846121862
$ws.Range("M20").Value = 6.8015495
$ws.Range("N20").Value = 13.603099
$ws.Range("O20").Value = 0.3428756056708687
$ws.Range("P20").Value = 0.2580802061075034
$ws.Range("Q20").Value = 526.937274202539
$ws.Range("R20").Value = 2107.749096810156
$ws.Range("S20").Value = 0.04034252027013781
$ws.Range("T20").Value = 0.0210700624922029

$ws.Range("G21").Value = 77.473122
$ws.Range("H21").Value = 154.946244
$ws.Range("I21").Value = 0.1176593481802354
$ws.Range("J21").Value = 0.08164152846121862
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.04027333333333333
$ws.Range("N21").Value = 0.12082
$ws.Range("O21").Value = 0.002030234957350741
$ws.Range("P21").Value = 0.002292216685470609
$ws.Range("Q21").Value = 3.12010086668
$ws.Range("R21").Value = 18.72060520008
$ws.Range("S21").Value = 0.0002388761217346161
$ws.Range("T21").Value = 0.0001871400737661289

$ws.Range("G22").Value = 145.7496183333334
$ws.Range("H22").Value = 437.248855
$ws.Range("I22").Value = 0.2213516720110761
$ws.Range("J22").Value = 0.2303874164262914
$ws.Range("M22").Value = 12.431794
$ws.Range("N22").Value = 37.295382
$ws.Range("O22").Value = 0.6267040910788743
$ws.Range("P22").Value = 0.7075740515758999
$ws.Range("Q22").Value = 1811.929230698624
$ws.Range("R22").Value = 16307.36307628761
$ws.Range("S22").Value = 0.1387219984164905
$ws.Range("T22").Value = 0.163016157672855

$ws.Range("G23").Value = 145.7496183333334
$ws.Range("H23").Value = 437.248855
$ws.Range("I23").Value = 0.2213516720110761
$ws.Range("J23").Value = 0.2303874164262914
$ws.Range("O23").Value = 0.0264162940991436
$ws.Range("P23").Value = 0.0298250554119953
$ws.Range("Q23").Value = 76.37488908452944
$ws.Range("R23").Value = 687.3740017607651
$ws.Range("S23").Value = 0.005847290867181758
$ws.Range("T23").Value = 0.006871317461140578

$ws.Range("G24").Value = 145.7496183333334
$ws.Range("H24").Value = 437.248855
$ws.Range("I24").Value = 0.2213516720110761
$ws.Range("J24").Value = 0.2303874164262914
$ws.Range("M24").Value = 0.03915333333333333
$ws.Range("N24").Value = 0.11746
$ws.Range("O24").Value = 0.001973774193762771
$ws.Range("P24").Value = 0.002228470219130754
$ws.Range("Q24").Value = 5.706583389811112
$ws.Range("R24").Value = 51.3592505083
$ws.Range("S24").Value = 0.000436898217961703
$ws.Range("T24").Value = 0.0005134114963684659

$ws.Range("G25").Value = 145.7496183333334
$ws.Range("H25").Value = 437.248855
$ws.Range("I25").Value = 0.2213516720110761
$ws.Range("J25").Value = 0.2303874164262914
$ws.Range("M25").Value = 6.8015495
$ws.Range("N25").Value = 13.603099
$ws.Range("O25").Value = 0.3428756056708687
$ws.Range("P25").Value = 0.2580802061075034
$ws.Range("Q25").Value = 991.3232437002744
$ws.Range("R25").Value = 5947.939462201646
$ws.Range("S25").Value = 0.07589608860705718
$ws.Range("T25").Value = 0.05945843191587249

$ws.Range("G26").Value = 145.7496183333334
$ws.Range("H26").Value = 437.248855
$ws.Range("I26").Value = 0.2213516720110761
$ws.Range("J26").Value = 0.2303874164262914
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0.3333333333333333
$ws.Range("M26").Value = 0.04027333333333333
$ws.Range("N26").Value = 0.12082
$ws.Range("O26").Value = 0.002030234957350741
$ws.Range("P26").Value = 0.002292216685470609
$ws.Range("Q26").Value = 5.869822962344446
$ws.Range("R26").Value = 52.8284066611
$ws.Range("S26").Value = 0.0004493959023849223
$ws.Range("T26").Value = 0.0005280978800548107
